# New weekly price report: insert a new data row at row 101, pushing the
# existing rows 101-143 down to 102-144, then populate the new row 101
# with the latest report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 101, shifting everything below it down.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new weekly record.
$ws.Range("A101").Value = 1
$ws.Range("B101").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C101").Value = "Arica y Parinacota"
$ws.Range("D101").Value = 44755
$ws.Range("E101").Value = 15
$ws.Range("F101").Value = "Fruta"
$ws.Range("G101").Value = 100108
$ws.Range("H101").Value = "Tropicales y subtropicales"
$ws.Range("I101").Value = 100108002
$ws.Range("J101").Value = "Mango"
$ws.Range("K101").Value = "Sin especificar"
$ws.Range("L101").Value = "Especial"
$ws.Range("M101").Value = 250
$ws.Range("N101").Value = 7000
$ws.Range("O101").Value = 7500
$ws.Range("P101").Value = 7250
$ws.Range("Q101").Value = "$/bandeja 4 kilos"
$ws.Range("R101").Value = "Brasil"
$ws.Range("S101").Value = 1812
$ws.Range("T101").Value = 4
